$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2033898305084746
$ws.Range("C2").Value = 0.5322033898305085
$ws.Range("J2").Value = 0.01016949152542373
$ws.Range("O2").Value = 0.003389830508474576
$ws.Range("P2").Value = 0.1525423728813559
$ws.Range("S2").Value = 0.09830508474576272
$ws.Range("B3").Value = 0.006289308176100629
$ws.Range("C3").Value = 0.03144654088050314
$ws.Range("J3").Value = 0.0440251572327044
$ws.Range("P3").Value = 0.7044025157232704
$ws.Range("S3").Value = 0.2138364779874214
$ws.Range("J4").Value = 0.02
$ws.Range("P4").Value = 0.74
$ws.Range("S4").Value = 0.24
$ws.Range("B6").Value = 0.05660377358490566
$ws.Range("D6").Value = 0.01509433962264151
$ws.Range("E6").Value = 0.007547169811320755
$ws.Range("F6").Value = 0.06037735849056604
$ws.Range("J6").Value = 0.3547169811320754
$ws.Range("O6").Value = 0.01509433962264151
$ws.Range("Q6").Value = 0.1245283018867925
$ws.Range("R6").Value = 0.06792452830188679
$ws.Range("S6").Value = 0.2981132075471698
$ws.Range("B7").Value = 0.1085972850678733
$ws.Range("D7").Value = 0.02714932126696833
$ws.Range("F7").Value = 0.07239819004524888
$ws.Range("J7").Value = 0.1266968325791855
$ws.Range("O7").Value = 0.01357466063348416
$ws.Range("Q7").Value = 0.1493212669683258
$ws.Range("R7").Value = 0.09954751131221719
$ws.Range("S7").Value = 0.4027149321266968
$ws.Range("B8").Value = 0.09580838323353294
$ws.Range("D8").Value = 0.01197604790419162
$ws.Range("F8").Value = 0.03393213572854292
$ws.Range("J8").Value = 0.1756487025948104
$ws.Range("O8").Value = 0.02594810379241517
$ws.Range("Q8").Value = 0.1377245508982036
$ws.Range("R8").Value = 0.1037924151696607
$ws.Range("S8").Value = 0.4151696606786427
$ws.Range("B9").Value = 0.09326424870466321
$ws.Range("D9").Value = 0.02590673575129534
$ws.Range("F9").Value = 0.07253886010362694
$ws.Range("J9").Value = 0.1139896373056995
$ws.Range("O9").Value = 0.0155440414507772
$ws.Range("Q9").Value = 0.1761658031088083
$ws.Range("R9").Value = 0.08808290155440414
$ws.Range("S9").Value = 0.4145077720207254
$ws.Range("B10").Value = 0.08779443254817987
$ws.Range("D10").Value = 0.02498215560314061
$ws.Range("F10").Value = 0.07066381156316917
$ws.Range("J10").Value = 0.1384725196288366
$ws.Range("O10").Value = 0.02498215560314061
$ws.Range("Q10").Value = 0.1862955032119914
$ws.Range("R10").Value = 0.092790863668808
$ws.Range("S10").Value = 0.3740185581727338
$ws.Range("G11").Value = 0.1270358306188925
$ws.Range("J11").Value = 0.1107491856677524
$ws.Range("K11").Value = 0.1954397394136808
$ws.Range("L11").Value = 0.5570032573289903
$ws.Range("S11").Value = 0.009771986970684038
$ws.Range("G12").Value = 0.7783783783783784
$ws.Range("J12").Value = 0.1675675675675676
$ws.Range("K12").Value = 0.005405405405405406
$ws.Range("L12").Value = 0.03243243243243243
$ws.Range("S12").Value = 0.01621621621621622
$ws.Range("G13").Value = 0.6909090909090909
$ws.Range("J13").Value = 0.2727272727272727
$ws.Range("S13").Value = 0.03636363636363636
$ws.Range("F15").Value = 0.03389830508474576
$ws.Range("H15").Value = 0.1398305084745763
$ws.Range("I15").Value = 0.08898305084745763
$ws.Range("J15").Value = 0.3177966101694915
$ws.Range("K15").Value = 0.05508474576271186
$ws.Range("M15").Value = 0.01694915254237288
$ws.Range("O15").Value = 0.0635593220338983
$ws.Range("S15").Value = 0.2838983050847458
$ws.Range("F16").Value = 0.02688172043010753
$ws.Range("H16").Value = 0.1827956989247312
$ws.Range("I16").Value = 0.06989247311827956
$ws.Range("J16").Value = 0.3548387096774194
$ws.Range("K16").Value = 0.1236559139784946
$ws.Range("M16").Value = 0.01612903225806452
$ws.Range("O16").Value = 0.07526881720430108
$ws.Range("S16").Value = 0.1505376344086022
$ws.Range("F17").Value = 0.03764705882352941
$ws.Range("H17").Value = 0.2023529411764706
$ws.Range("I17").Value = 0.08
$ws.Range("J17").Value = 0.4235294117647059
$ws.Range("K17").Value = 0.08941176470588236
$ws.Range("M17").Value = 0.01647058823529412
$ws.Range("N17").Value = 0.004705882352941176
$ws.Range("O17").Value = 0.04941176470588235
$ws.Range("S17").Value = 0.09647058823529411
$ws.Range("F18").Value = 0.02966101694915254
$ws.Range("H18").Value = 0.1822033898305085
$ws.Range("I18").Value = 0.08898305084745763
$ws.Range("J18").Value = 0.4533898305084746
$ws.Range("K18").Value = 0.1101694915254237
$ws.Range("M18").Value = 0.02542372881355932
$ws.Range("O18").Value = 0.0423728813559322
$ws.Range("S18").Value = 0.06779661016949153
$ws.Range("F19").Value = 0.0270473328324568
$ws.Range("H19").Value = 0.2283996994740796
$ws.Range("I19").Value = 0.07663410969196092
$ws.Range("J19").Value = 0.3613824192336589
$ws.Range("K19").Value = 0.1074380165289256
$ws.Range("M19").Value = 0.02629601803155522
$ws.Range("N19").Value = 0.001502629601803156
$ws.Range("O19").Value = 0.07062359128474831
$ws.Range("S19").Value = 0.1006761833208114
